$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 107
$ws.Range("I33").Value = 73.40000000000001
$ws.Range("K33").Value = 73.40000000000001
$ws.Range("M33").Value = 155.6
$ws.Range("H40").Value = 7175.0835
$ws.Range("I40").Value = 5872.2856
$ws.Range("J40").Value = 8999
$ws.Range("K40").Value = 5872.2856
$ws.Range("L40").Value = 8999
$ws.Range("M40").Value = -5697.2856
$ws.Range("N40").Value = -9349
$ws.Range("H64").Value = 3799.5
$ws.Range("I64").Value = 3799
$ws.Range("J64").Value = 3800
$ws.Range("K64").Value = 3799
$ws.Range("L64").Value = 3800
$ws.Range("M64").Value = -3551
$ws.Range("N64").Value = -4296
$ws.Range("H67").Value = 3799.5
$ws.Range("I67").Value = 3799
$ws.Range("J67").Value = 3800
$ws.Range("K67").Value = 3799
$ws.Range("L67").Value = 3800
$ws.Range("M67").Value = -2941
$ws.Range("N67").Value = -5516
$ws.Range("H82").Value = 715.3333
$ws.Range("I82").Value = 715.3333
$ws.Range("K82").Value = 2145.9999
$ws.Range("M82").Value = -1739.9999
$ws.Range("H85").Value = 715.3333
$ws.Range("I85").Value = 715.3333
$ws.Range("K85").Value = 2145.9999
$ws.Range("M85").Value = -741.9998999999998
$ws.Range("H111").Value = 1767.2903
$ws.Range("I111").Value = 921.5
$ws.Range("J111").Value = 1892.5927
$ws.Range("K111").Value = 2764.5
$ws.Range("L111").Value = 5677.7781
$ws.Range("M111").Value = 302.5
$ws.Range("N111").Value = -11811.7781
$ws.Range("H132").Value = 1874
$ws.Range("I132").Value = 1368.2727
$ws.Range("K132").Value = 4104.8181
$ws.Range("M132").Value = -1574.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H37").Value = 27777.666
$ws.Range("I37").Value = 17500
$ws.Range("J37").Value = 35999.8
$ws.Range("K37").Value = 17500
$ws.Range("L37").Value = 35999.8
$ws.Range("M37").Value = -17227
$ws.Range("N37").Value = -36545.8
$ws.Range("H102").Value = 12506089
$ws.Range("I102").Value = 25002398
$ws.Range("K102").Value = 25002398
$ws.Range("M102").Value = -25000776
$ws.Range("H122").Value = 2774
$ws.Range("I122").Value = 2774
$ws.Range("K122").Value = 8322
$ws.Range("M122").Value = -5872

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2911.8096
$ws.Range("I86").Value = 1331.7273
$ws.Range("J86").Value = 4649.9
$ws.Range("K86").Value = 1331.7273
$ws.Range("L86").Value = 4649.9
$ws.Range("M86").Value = -208.7273
$ws.Range("N86").Value = -6895.9
$ws.Range("H88").Value = 37528.43
$ws.Range("J88").Value = 37528.43
$ws.Range("L88").Value = 37528.43
$ws.Range("N88").Value = -38340.43
$ws.Range("H89").Value = 2911.8096
$ws.Range("I89").Value = 1331.7273
$ws.Range("J89").Value = 4649.9
$ws.Range("K89").Value = 6658.636500000001
$ws.Range("L89").Value = 23249.5
$ws.Range("M89").Value = -1042.636500000001
$ws.Range("N89").Value = -34481.5
$ws.Range("H91").Value = 37528.43
$ws.Range("J91").Value = 37528.43
$ws.Range("L91").Value = 37528.43
$ws.Range("N91").Value = -40336.43
$ws.Range("H107").Value = 62506428
$ws.Range("I107").Value = 166667140
$ws.Range("J107").Value = 9998
$ws.Range("K107").Value = 166667140
$ws.Range("L107").Value = 9998
$ws.Range("M107").Value = -166665220
$ws.Range("N107").Value = -13838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 781.5
$ws.Range("I22").Value = 448.54544
$ws.Range("K22").Value = 448.54544
$ws.Range("M22").Value = -98.54543999999999
$ws.Range("H58").Value = 3446.25
$ws.Range("I58").Value = 1714.2727
$ws.Range("K58").Value = 1714.2727
$ws.Range("M58").Value = -1511.2727
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = ""
$ws.Range("H107").Value = 1125.6428
$ws.Range("I107").Value = 465.1111
$ws.Range("J107").Value = 2314.6
$ws.Range("K107").Value = 465.1111
$ws.Range("L107").Value = 2314.6
$ws.Range("M107").Value = 1454.8889
$ws.Range("N107").Value = -6154.6
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 1819.1333
$ws.Range("I132").Value = 1819.1333
$ws.Range("K132").Value = 5457.3999
$ws.Range("M132").Value = -2927.3999
$ws.Range("H134").Value = 3171.0833
$ws.Range("I134").Value = 2549
$ws.Range("K134").Value = 7647
$ws.Range("M134").Value = -5112
$ws.Range("H136").Value = 3446.25
$ws.Range("I136").Value = 1714.2727
$ws.Range("K136").Value = 5142.8181
$ws.Range("M136").Value = -2592.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1802.6666
$ws.Range("I34").Value = 185.4
$ws.Range("J34").Value = 2957.8572
$ws.Range("K34").Value = 556.2
$ws.Range("L34").Value = 8873.571599999999
$ws.Range("M34").Value = -472.2
$ws.Range("N34").Value = -9041.571599999999
$ws.Range("H40").Value = 162.55556
$ws.Range("J40").Value = 330
$ws.Range("L40").Value = 1320
$ws.Range("N40").Value = -1458
$ws.Range("H54").Value = 1425
$ws.Range("I54").Value = 1300
$ws.Range("J54").Value = 1487.5
$ws.Range("K54").Value = 3900
$ws.Range("L54").Value = 4462.5
$ws.Range("M54").Value = -3341
$ws.Range("N54").Value = -5580.5
$ws.Range("H103").Value = 2410.5
$ws.Range("J103").Value = 2410.5
$ws.Range("L103").Value = 7231.5
$ws.Range("N103").Value = -8989.5
$ws.Range("H113").Value = 1498.3334
$ws.Range("I113").Value = 1495
$ws.Range("K113").Value = 4485
$ws.Range("M113").Value = -2315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 500
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 500
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H97").Value = 844.5
$ws.Range("I97").Value = 189
$ws.Range("K97").Value = 189
$ws.Range("M97").Value = 307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1087.7693
$ws.Range("I16").Value = 1087.7693
$ws.Range("K16").Value = 1087.7693
$ws.Range("M16").Value = -917.7692999999999
$ws.Range("H46").Value = 6933.3335
$ws.Range("I46").Value = 5333.3335
$ws.Range("J46").Value = 8000
$ws.Range("K46").Value = 5333.3335
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = -5145.3335
$ws.Range("N46").Value = -8376
$ws.Range("H82").Value = 2866.375
$ws.Range("I82").Value = 512.55554
$ws.Range("J82").Value = 5892.7144
$ws.Range("K82").Value = 512.55554
$ws.Range("L82").Value = 5892.7144
$ws.Range("M82").Value = -151.55554
$ws.Range("N82").Value = -6614.7144
$ws.Range("H85").Value = 2866.375
$ws.Range("I85").Value = 512.55554
$ws.Range("J85").Value = 5892.7144
$ws.Range("K85").Value = 512.55554
$ws.Range("L85").Value = 5892.7144
$ws.Range("M85").Value = 735.44446
$ws.Range("N85").Value = -8388.714400000001
$ws.Range("H122").Value = 3228.7144
$ws.Range("I122").Value = 4118.4
$ws.Range("J122").Value = 1004.5
$ws.Range("K122").Value = 12355.2
$ws.Range("L122").Value = 3013.5
$ws.Range("M122").Value = -9905.199999999999
$ws.Range("N122").Value = -7913.5
$ws.Range("H132").Value = 2967
$ws.Range("I132").Value = 2967
$ws.Range("K132").Value = 8901
$ws.Range("M132").Value = -6371

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11500
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 11500
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -66240
$ws.Range("H107").Value = 33334218
$ws.Range("I107").Value = 41667370
$ws.Range("J107").Value = 1625
$ws.Range("K107").Value = 125002110
$ws.Range("L107").Value = 4875
$ws.Range("M107").Value = -125000190
$ws.Range("N107").Value = -8715
$ws.Range("H122").Value = 2124.4285
$ws.Range("J122").Value = 2015.3334
$ws.Range("L122").Value = 6046.0002
$ws.Range("N122").Value = -10946.0002
